$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
    return $rng.Find.Found
}

# --- 1. "Date of Next review: Spring Term 2020" -> "...Autumn Term 2021" (must run before the
#        "Summer"->"Spring" replacement below, while "Spring" is still unique) ---
Replace-Text "Date of Next review: Spring Term 2020" "Date of Next review: Autumn Term 2021" | Out-Null

# --- 2. "Written: Summer Term 2018" -> "Written: Spring Term 2020" ---
Replace-Text "Written: Summer Term 2018" "Written: Spring Term 2020" | Out-Null

# --- 3. "Approved by: ... 31st July 2018" -> "...18th February 2020" ---
Replace-Text "31st July 2018" "18th February 2020" | Out-Null

Write-Output "stage1 done"

# --- 4. Paragraphs where Word's proofing "grammar" marks (in order to / public sector /
#        In particular, the / in the area of / on account of / in the course of) were
#        accepted/cleared, leaving the surrounding runs merged together. The visible text
#        is unchanged, so we replace each run of text with itself; doing so via Find &
#        Replace causes the host to rebuild the run (dropping the now-stale proofErr marks). ---

$q1 = [char]0x201C
$q2 = [char]0x201D
$apos = [char]0x2019

$t1 = "Our admissions arrangements are fair and transparent, and do not discriminate on race, disability, sexual orientation or socio-economic factors. Exclusions will always be based on the school" + $apos + "s Behaviour Policy. We will closely monitor exclusions in respect of equality in order to avoid any potential adverse impact.   "
Replace-Text $t1 $t1 | Out-Null

$t2 = "The Equality Act (2010) introduced a single equality duty for all public sector organisations including schools: known as the " + $q1 + "public sector equality duty" + $q2 + ".  The public sector duty requires all schools to show how they are meeting the aims of the Equality Act by giving due regard to the need to: "
Replace-Text $t2 $t2 | Out-Null

$t3 = "This Equality Policy and Plan is our response to demonstrate a) that we comply with the duty to have due regard for the three aims above and b) specific and measurable objectives which will be pursued over the coming years to achieve the three aims. In particular, the action plan at the end of this Equality Policy and Plan outlines the actions "
Replace-Text $t3 $t3 | Out-Null

$t4 = "Legislation makes discrimination unlawful in the area of goods, facilities and services on grounds of sexual orientation. For schools this means admissions, benefits and services for students and treatment of students.   "
Replace-Text $t4 $t4 | Out-Null

$t5 = "Harassment on account of race, gender, disability or sexual orientation is unacceptable and is not tolerated within the school environment.   All staff are expected to deal with any discriminatory incidents that may occur. They are expected to know how to identify and challenge prejudice and stereotyping; and to support the full range of diverse needs according to a studen"
Replace-Text $t5 $t5 | Out-Null

$t6 = [char]0x2022 + " discriminatory comments in the course of discussion; "
Replace-Text $t6 $t6 | Out-Null

Write-Output "stage2 done"
